$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.526.44"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.650.17"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "302.23"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").Value = "0.3597"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "51.01"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").Value = "0.08202"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").Value = "1.234"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "22.35"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.470"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "7.487"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "0.00001226"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "1.647.74"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "97.54"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "0.06999"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "6.784"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").Value = "17.57"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "12.66"
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").Value = "23.536.67"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "2.518"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "3.028"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Value = "153.66"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "5.226"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "133.97"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").Value = "1.829.82"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").Value = "7.108"
$ws.Range("E32").Value = "  +10.24%  "
$ws.Range("D33").Value = "2.246"
$ws.Range("E33").Value = "  +6.41%  "
$ws.Range("D34").Value = "12.05"
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").Value = "1.061"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "0.02795"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").Value = "0.2496"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "6.092"
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("D39").Value = "0.08764"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "0.06997"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "13.12"
$ws.Range("E41").Value = "  +10.27%  "
$ws.Range("D42").Value = "0.6991"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").Value = "1.337"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").Value = "15.97"
$ws.Range("E44").Value = "  +4.62%  "
$ws.Range("D45").Value = "0.6519"
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "2.302"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").Value = "3.955"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07880"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +1.28%  "
